$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): relabel the three columns ---
# A1: "Email" -> "EmpleadoNo"
# B1: "Nuevo Email" -> "Empresa Id"
# C1: "Empresa Id" -> "Email" (and this cell loses its special styling)
$ws.Range("A1").Value = "EmpleadoNo"
$ws.Range("B1").Value = "Empresa Id"
$ws.Range("C1").Value = "Email"

# C1 reverts to the plain default look (no custom font/alignment anymore)
$ws.Range("C1").ClearFormats()

# --- Data rows (2 & 3) formatting tweaks ---
# Column A (EmpleadoNo) data cells: centered vertically, keep wrap, default font
$ws.Range("A2:A3").ClearFormats()
$ws.Range("A2:A3").WrapText = $true
$ws.Range("A2:A3").VerticalAlignment = -4108

# Column C (Email) data cells: styled as hyperlinks (email addresses), centered + wrapped
$ws.Range("C2:C3").Style = "Hyperlink"
$ws.Range("C2:C3").WrapText = $true
$ws.Range("C2:C3").VerticalAlignment = -4108

# --- Column widths to match the refreshed template layout ---
$ws.Columns.Item(1).ColumnWidth = 28.16666667
$ws.Columns.Item(2).ColumnWidth = 25.25
$ws.Columns.Item(3).ColumnWidth = 24.16666667

# --- Row heights explicit (matches default 15.75 but now set per-row) ---
$ws.Rows.Item(1).RowHeight = 15.75
$ws.Rows.Item(2).RowHeight = 15.75
$ws.Rows.Item(3).RowHeight = 15.75

# --- Selection / active cell state ---
$ws.Range("B14").Select()
